$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 1189.1
$ws_ALC.Range("I28").Value = 866.38464
$ws_ALC.Range("K28").Value = 866.38464
$ws_ALC.Range("M28").Value = -381.38464
$ws_ALC.Range("H69").Value = 57908.082
$ws_ALC.Range("J69").Value = 82362.125
$ws_ALC.Range("L69").Value = 247086.375
$ws_ALC.Range("N69").Value = -248834.375
$ws_ALC.Range("H72").Value = 57908.082
$ws_ALC.Range("J72").Value = 82362.125
$ws_ALC.Range("L72").Value = 741259.125
$ws_ALC.Range("N72").Value = -749995.125
$ws_ALC.Range("H92").Value = 975.7917
$ws_ALC.Range("I92").Value = 849
$ws_ALC.Range("J92").Value = 1187.1111
$ws_ALC.Range("K92").Value = 849
$ws_ALC.Range("L92").Value = 1187.1111
$ws_ALC.Range("M92").Value = 399
$ws_ALC.Range("N92").Value = -3683.1111
$ws_ALC.Range("H107").Value = 249
$ws_ALC.Range("I107").Value = 249
$ws_ALC.Range("J107").Value = 0
$ws_ALC.Range("K107").Value = 249
$ws_ALC.Range("L107").Value = 0
$ws_ALC.Range("M107").Value = 1671
$ws_ALC.Range("N107").ClearContents()
$ws_ALC.Range("H112").Value = 1982.2693
$ws_ALC.Range("I112").Value = 1949.5
$ws_ALC.Range("J112").Value = 1985
$ws_ALC.Range("K112").Value = 5848.5
$ws_ALC.Range("L112").Value = 5955
$ws_ALC.Range("M112").Value = -4740.5
$ws_ALC.Range("N112").Value = -8171
$ws_ALC.Range("H131").Value = 15182.5
$ws_ALC.Range("I131").Value = 14023.75
$ws_ALC.Range("K131").Value = 42071.25
$ws_ALC.Range("M131").Value = -37031.25
$ws_ALC.Range("H132").Value = 12689.945
$ws_ALC.Range("I132").Value = 2961.0444
$ws_ALC.Range("K132").Value = 8883.1332
$ws_ALC.Range("M132").Value = -6353.1332
$ws_ALC.Range("H137").Value = 5807.4688
$ws_ALC.Range("I137").Value = 5973.92
$ws_ALC.Range("J137").Value = 5213
$ws_ALC.Range("K137").Value = 17921.76
$ws_ALC.Range("L137").Value = 15639
$ws_ALC.Range("M137").Value = -15371.76
$ws_ALC.Range("N137").Value = -20739
$ws_ALC.Range("H138").Value = 3229.6296
$ws_ALC.Range("J138").Value = 4082.2104
$ws_ALC.Range("L138").Value = 12246.6312
$ws_ALC.Range("N138").Value = -22526.6312

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H6").Value = 0
$ws_ARM.Range("I6").Value = 0
$ws_ARM.Range("K6").Value = 0
$ws_ARM.Range("M6").ClearContents()
$ws_ARM.Range("H32").Value = 17248022
$ws_ARM.Range("I32").Value = 20415638
$ws_ARM.Range("K32").Value = 20415638
$ws_ARM.Range("M32").Value = -20415351
$ws_ARM.Range("H102").Value = 1933.5
$ws_ARM.Range("J102").Value = 1613.5
$ws_ARM.Range("L102").Value = 1613.5
$ws_ARM.Range("N102").Value = -4857.5
$ws_ARM.Range("H132").Value = 1032.9517
$ws_ARM.Range("I132").Value = 807.2245
$ws_ARM.Range("J132").Value = 1883.7693
$ws_ARM.Range("K132").Value = 2421.6735
$ws_ARM.Range("L132").Value = 5651.3079
$ws_ARM.Range("M132").Value = 108.3265000000001
$ws_ARM.Range("N132").Value = -10711.3079

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 2870.5715
$ws_BSM.Range("I94").Value = 3384.6365
$ws_BSM.Range("J94").Value = 985.6667
$ws_BSM.Range("K94").Value = 3384.6365
$ws_BSM.Range("L94").Value = 985.6667
$ws_BSM.Range("M94").Value = -2933.6365
$ws_BSM.Range("N94").Value = -1887.6667
$ws_BSM.Range("H107").Value = 3555.5
$ws_BSM.Range("I107").Value = 3856.5
$ws_BSM.Range("K107").Value = 3856.5
$ws_BSM.Range("M107").Value = -1936.5
$ws_BSM.Range("H132").Value = 69522.38
$ws_BSM.Range("J132").Value = 69522.38
$ws_BSM.Range("L132").Value = 69522.38
$ws_BSM.Range("N132").Value = -79642.38
$ws_BSM.Range("H134").Value = 2074.739
$ws_BSM.Range("I134").Value = 1582.4524
$ws_BSM.Range("K134").Value = 4747.357199999999
$ws_BSM.Range("M134").Value = -2212.357199999999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H122").Value = 302753
$ws_CRP.Range("I122").Value = 537959.9
$ws_CRP.Range("K122").Value = 1613879.7
$ws_CRP.Range("M122").Value = -1611429.7
$ws_CRP.Range("H132").Value = 2026.9459
$ws_CRP.Range("I132").Value = 1853.1471
$ws_CRP.Range("K132").Value = 5559.4413
$ws_CRP.Range("M132").Value = -3029.4413
$ws_CRP.Range("H134").Value = 2890.5
$ws_CRP.Range("I134").Value = 2575.05
$ws_CRP.Range("K134").Value = 7725.150000000001
$ws_CRP.Range("M134").Value = -5190.150000000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 78.90000000000001
$ws_CUL.Range("I2").Value = 63.6
$ws_CUL.Range("J2").Value = 94.2
$ws_CUL.Range("K2").Value = 381.6
$ws_CUL.Range("L2").Value = 565.2
$ws_CUL.Range("M2").Value = -268.6
$ws_CUL.Range("N2").Value = -791.2
$ws_CUL.Range("H34").Value = 101.72727
$ws_CUL.Range("I34").Value = 101.9
$ws_CUL.Range("K34").Value = 305.7
$ws_CUL.Range("M34").Value = -221.7
$ws_CUL.Range("H81").Value = 4332.6665
$ws_CUL.Range("I81").Value = 3999
$ws_CUL.Range("K81").Value = 11997
$ws_CUL.Range("M81").Value = -10874
$ws_CUL.Range("H84").Value = 4332.6665
$ws_CUL.Range("I84").Value = 3999
$ws_CUL.Range("K84").Value = 35991
$ws_CUL.Range("M84").Value = -30375
$ws_CUL.Range("H107").Value = 269.2857
$ws_CUL.Range("I107").Value = 239.4
$ws_CUL.Range("J107").Value = 344
$ws_CUL.Range("K107").Value = 718.2
$ws_CUL.Range("L107").Value = 1032
$ws_CUL.Range("M107").Value = 1201.8
$ws_CUL.Range("N107").Value = -4872
$ws_CUL.Range("H113").Value = 1765.9
$ws_CUL.Range("I113").Value = 3725
$ws_CUL.Range("J113").Value = 1548.2222
$ws_CUL.Range("K113").Value = 11175
$ws_CUL.Range("L113").Value = 4644.6666
$ws_CUL.Range("M113").Value = -9005
$ws_CUL.Range("N113").Value = -8984.6666
$ws_CUL.Range("H131").Value = 3824.7222
$ws_CUL.Range("I131").Value = 3041.5
$ws_CUL.Range("J131").Value = 4216.3335
$ws_CUL.Range("K131").Value = 9124.5
$ws_CUL.Range("L131").Value = 12649.0005
$ws_CUL.Range("M131").Value = -4084.5
$ws_CUL.Range("N131").Value = -22729.0005
$ws_CUL.Range("H139").Value = 76926696
$ws_CUL.Range("I139").Value = 76926696
$ws_CUL.Range("K139").Value = 230780088
$ws_CUL.Range("M139").Value = -230774948

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 9560.866
$ws_LTW.Range("I40").Value = 12150.8
$ws_LTW.Range("K40").Value = 12150.8
$ws_LTW.Range("M40").Value = -12014.8
$ws_LTW.Range("H100").Value = 61850.1
$ws_LTW.Range("I100").Value = 162450.42
$ws_LTW.Range("J100").Value = 7680.6924
$ws_LTW.Range("K100").Value = 162450.42
$ws_LTW.Range("L100").Value = 7680.6924
$ws_LTW.Range("M100").Value = -161909.42
$ws_LTW.Range("N100").Value = -8762.6924
$ws_LTW.Range("H122").Value = 4294.8887
$ws_LTW.Range("I122").Value = 1676.0416
$ws_LTW.Range("K122").Value = 5028.1248
$ws_LTW.Range("M122").Value = -2578.1248
$ws_LTW.Range("H132").Value = 4744.263
$ws_LTW.Range("I132").Value = 3233.4583
$ws_LTW.Range("K132").Value = 9700.374899999999
$ws_LTW.Range("M132").Value = -7170.374899999999

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H54").Value = 28165
$ws_WVR.Range("J54").Value = 28165
$ws_WVR.Range("L54").Value = 28165
$ws_WVR.Range("N54").Value = -29205
$ws_WVR.Range("H119").Value = 28500
$ws_WVR.Range("J119").Value = 28500
$ws_WVR.Range("L119").Value = 28500
$ws_WVR.Range("N119").Value = -38176
$ws_WVR.Range("H122").Value = 2357.205
$ws_WVR.Range("I122").Value = 2002.1515
$ws_WVR.Range("K122").Value = 6006.4545
$ws_WVR.Range("M122").Value = -3556.4545
$ws_WVR.Range("H132").Value = 1572.4872
$ws_WVR.Range("I132").Value = 1387.2162
$ws_WVR.Range("K132").Value = 4161.6486
$ws_WVR.Range("M132").Value = -1631.6486
